$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 72
    4  = 266
    5  = 330
    6  = 84
    7  = 132
    8  = 92
    9  = 79
    10 = 49
    11 = 215
    12 = 468
    13 = 320
    14 = 103
    15 = 66
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 8).Value = $updates[$row]
}
